$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week1")

# Update existing rows 4, 6, 7
$ws.Range("D4").Value = "N"
$ws.Range("E4").Value = "Rachel Fontaine"
$ws.Range("F4").Value = 150.0

$ws.Range("C6").Value = "SCH4U"
$ws.Range("D6").Value = "N"
$ws.Range("E6").Value = "Leonard Robillard"
$ws.Range("F6").Value = 136.0
$ws.Range("G6").Value = "Fletcher Donaldson"
$ws.Range("H6").Value = 120.0

$ws.Range("B7").Value = 2.0
$ws.Range("C7").Value = "SCH3U"
$ws.Range("E7").Value = "Dineth"
$ws.Range("F7").Value = 165.0

# Add new rows 8-11
$ws.Range("A8").Value = 1.0
$ws.Range("B8").Value = 4.0
$ws.Range("C8").Value = "SNC2P"
$ws.Range("D8").Value = "N"
$ws.Range("E8").Value = "Allison Petersen"
$ws.Range("F8").Value = 103.0
$ws.Range("G8").Value = "Fletcher Donaldson"
$ws.Range("H8").Value = 120.0

$ws.Range("A9").Value = 1.0
$ws.Range("B9").Value = 1.0
$ws.Range("C9").Value = "ENG1D"
$ws.Range("D9").Value = "N"
$ws.Range("E9").Value = "Miaomiao Yan"
$ws.Range("F9").Value = 142.0
$ws.Range("G9").Value = "Guadalupe Case"
$ws.Range("H9").Value = 124.0

$ws.Range("A10").Value = 1.0
$ws.Range("B10").Value = 2.0
$ws.Range("C10").Value = "ENG3U"
$ws.Range("D10").Value = "N"
$ws.Range("E10").Value = "Merlin Lambert"
$ws.Range("F10").Value = 141.0
$ws.Range("G10").Value = "Guadalupe Case"
$ws.Range("H10").Value = 124.0

$ws.Range("A11").Value = 1.0
$ws.Range("B11").Value = 3.0
$ws.Range("C11").Value = "ENG1D"
$ws.Range("D11").Value = "N"
$ws.Range("E11").Value = "Tammy Mcmillan"
$ws.Range("F11").Value = 159.0
$ws.Range("G11").Value = "Guadalupe Case"
$ws.Range("H11").Value = 124.0
